$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-103). The diff changes each of these from 45190 to 45192.
for ($r = 2; $r -le 103; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
